$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Redefining Security: The Future of Cyber Defense" "The Artful Science: Marrying Aesthetics and Logic in Mathematics"

# Author name / title
Replace-Text "Dr" "Prof"
Replace-Text " Evelyn Richards" " Laura Harper"

# Email
Replace-Text "erichards@scienceinstitute" "lauraharper835@baysidehs"
Replace-Text "org" "edu"

# Intro paragraph sentences
Replace-Text "The rapid evolution of cyberspace has ushered in a new era of challenges and opportunities in the realm of security" "Within the vast realm of human knowledge, mathematics stands uniquely as both an art and a science, forging an inseparable bond between aesthetics and logic"
Replace-Text " With the meteoric rise of interconnected devices and the seamless integration of technology into every facet of our lives, the threat landscape has exponentially expanded" " Its essence lies in the exploration of patterns, the seeking of regularities amidst the apparent chaos of our world"
Replace-Text " As individuals, organizations, and nations strive to navigate this dynamic terrain, the conventional approach to cybersecurity is no longer adequate" " Through its intricate web of symbols, formulas, and theorems, mathematics unravels the underlying order and beauty that permeates our universe"
Replace-Text " It is imperative to embrace innovative strategies, cutting-edge technologies, and a proactive mindset to defend against sophisticated threats and maintain digital resilience" " As we traverse this realm, we discover intricate patterns in nature, marveling at the Fibonacci sequence echoing in the spiral tendrils of a fern or the mesmerizing symmetry of a sunflower"

Replace-Text "In this essay, we delve into the future of cyber defense, exploring the paradigm shift from reactive measures to proactive preparedness" "Mathematics finds its canvas in architecture, where it transforms into the graceful arches of bridges and soaring spires of cathedrals, each structure a testament to the harmony between form and function"
Replace-Text " We examine the pivotal role of artificial intelligence and machine learning in enhancing threat detection and response capabilities, enabling systems to learn from past incidents and adapt in real time" " In music, mathematics guides the symphony of sounds, revealing the underlying rhythmic patterns and harmonic progressions that stir our emotions and ignite our imaginations"
Replace-Text " We highlight the significance of collaboration and information sharing among stakeholders, emphasizing the need for public-private partnerships and international cooperation to effectively address global cyber threats" " Even in the works of great artists, we can discern the subtle hand of mathematics, whether it manifests as the golden ratio in Leonardo da Vinci's masterpieces or the fractal geometry that echoes throughout Jackson Pollock's abstract expressionist paintings"

Replace-Text "Moreover, we discuss the imperative to raise cybersecurity awareness and educate users about emerging threats and best practices" "The elegance and simplicity of a mathematical solution, the aha moment when a complex puzzle finally yields to logical reasoning, is an experience akin to witnessing the unfolding of a beautiful melody or the revelation of a hidden truth"
Replace-Text " By empowering individuals to recognize and mitigate vulnerabilities, we create a safer cyber ecosystem" " It is in these moments of mathematical enlightenment that we glimpse the profound interconnectedness of all knowledge, the unity that underlies the diversity of human expression"
Replace-Text " Lastly, we emphasize the critical role of adopting a comprehensive and holistic approach to cyber defense, encompassing technical, legal, and policy frameworks, to ensure a truly secure and resilient digital infrastructure" " As we delve deeper into the world of mathematics, we cultivate a deeper appreciation for its elegance, its power, and its universality, recognizing it as both a tool and a language that transcends cultural and linguistic boundaries"

# Summary paragraph
Replace-Text "The future of cyber defense demands a transformative approach that encompasses a multifaceted and proactive strategy" "Mathematics, at its core, is both an art and a science"
Replace-Text " We must leverage advanced technologies, foster collaboration, cultivate a culture of awareness, and adopt a comprehensive framework to safeguard our increasingly interconnected world" " It weaves together aesthetics and logic in a mesmerizing tapestry of patterns and harmonies, revealing the underlying order and beauty that permeate our world"
Replace-Text " The journey towards a secured cyber future relies on the collective efforts of individuals, organizations, and governments to create a safer digital landscape for all" " In the elegance of its solutions and the profound interconnectedness of its concepts, we glimpse the unity that binds all knowledge together"

Write-Host "Done with simple replacements"
